$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.683.34"
$ws.Range("E2").Value = "  -1.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.419.42"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.05"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.32"
$ws.Range("E6").Value = "  -2.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.607"
$ws.Range("E7").Value = "  +3.60%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.422.24"
$ws.Range("E9").Value = "  -1.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.13"
$ws.Range("E10").Value = "  -2.66%  "
$ws.Range("E11").Value = "  -2.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.439"
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.007.18"
$ws.Range("E13").Value = "  -2.02%  "
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("E15").Value = "  -4.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.55"
$ws.Range("E16").Value = "  -3.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.647.41"
$ws.Range("E17").Value = "  -1.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.389.73"
$ws.Range("E18").Value = "  -2.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.34"
$ws.Range("E19").Value = "  -1.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.85"
$ws.Range("E20").Value = "  -3.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "380.14"
$ws.Range("E21").Value = "  -3.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.98"
$ws.Range("E22").Value = "  -3.67%  "
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.20"
$ws.Range("E25").Value = "  -2.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000119"
$ws.Range("E26").Value = "  -5.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.10"
$ws.Range("E27").Value = "  +5.49%  "
$ws.Range("E28").Value = "  -1.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.48"
$ws.Range("E30").Value = "  +1.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.18"
$ws.Range("E31").Value = "  -4.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.01"
$ws.Range("E32").Value = "  -2.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.18"
$ws.Range("E33").Value = "  -2.48%  "
$ws.Range("E34").Value = "  -1.18%  "
$ws.Range("E35").Value = "  +2.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.31"
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.92"
$ws.Range("E37").Value = "  -3.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0756"
$ws.Range("E38").Value = "  -2.52%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.885.73"
$ws.Range("E39").Value = "  -5.69%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.79"
$ws.Range("E40").Value = "  +2.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.42"
$ws.Range("E41").Value = "  -3.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.59"
$ws.Range("E42").Value = "  +1.16%  "
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("E44").Value = "  -2.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.770"
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.71"
$ws.Range("E46").Value = "  -1.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "320.93"
$ws.Range("E47").Value = "  +2.70%  "
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.22"
$ws.Range("E48").Value = "  -1.65%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.07"
$ws.Range("E49").Value = "  -5.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.107"
$ws.Range("E50").Value = "  +0.39%  "
